# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (rId1) and "全部类型" (rId4) sheets, per the latest data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — row numbers per dataset
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 1336
$wsExpo.Range("F11").Value = 4776
$wsExpo.Range("F21").Value = 3698
$wsExpo.Range("F22").Value = 529
$wsExpo.Range("F34").Value = 832

# Sheet "全部类型" (All types) — same events, one row lower after row 34
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1336
$wsAll.Range("F11").Value = 4776
$wsAll.Range("F21").Value = 3698
$wsAll.Range("F22").Value = 529
$wsAll.Range("F35").Value = 832
